$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cryptocurrency price / volume figures (and the Algorand/ARBITRUM row swap)

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.227.69'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.05%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.912.41'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.26%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.25%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.75'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.99%  '

# Row 7
$ws.Range("E7").Value = '  +0.50%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3932'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.27%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09367'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.63%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.141'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.25%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.96'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.84%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.408'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.05%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.93'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.29%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.910.98'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.56%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.328'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.12%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9998'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.28%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001123'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.55%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.49'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.40%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06616'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.16%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '18.00'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.04%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9998'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.18%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.234'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.76%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.265.55'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.15%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.52'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.14%  '

# Row 25
$ws.Range("E25").Value = '  +1.32%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.602'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.17%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.125.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.79%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.14'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.24%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '158.03'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.45%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.30'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.24%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.108'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.03%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1076'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.01%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.656'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.48%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.613'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.28%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.732'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.48%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06685'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.24%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02428'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.36%  '

# Row 38
$ws.Range("B38").Value = 'ARBITRUM'
$ws.Range("C38").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.248'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.47%  '

# Row 39
$ws.Range("B39").Value = 'Algorand'
$ws.Range("C39").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2211'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.86%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.284'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +8.40%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6530'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.74%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.55'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.37%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.023'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.40%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9994'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.16%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6121'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.17%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.34'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.55%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.728'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.91%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.287'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.05%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.025'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.12%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '123.18'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.61%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.190'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.55%  '
